$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 821.2308
$ws.Range("I92").Value = 473
$ws.Range("J92").Value = 5000
$ws.Range("K92").Value = 473
$ws.Range("L92").Value = 5000
$ws.Range("M92").Value = 775
$ws.Range("N92").Value = -7496
$ws.Range("H99").Value = 499.5
$ws.Range("I99").Value = 499.5
$ws.Range("K99").Value = 1498.5
$ws.Range("M99").Value = -0.5
$ws.Range("H100").Value = 2587.9092
$ws.Range("J100").Value = 3593.4
$ws.Range("L100").Value = 3593.4
$ws.Range("N100").Value = -4675.4
$ws.Range("H113").Value = 31257398
$ws.Range("I113").Value = 3130.1428
$ws.Range("J113").Value = 55566270
$ws.Range("K113").Value = 3130.1428
$ws.Range("L113").Value = 55566270
$ws.Range("M113").Value = 123.8571999999999
$ws.Range("N113").Value = -55572778
$ws.Range("H132").Value = 1342.3103
$ws.Range("I132").Value = 1557.15
$ws.Range("J132").Value = 864.8889
$ws.Range("K132").Value = 4671.450000000001
$ws.Range("L132").Value = 2594.6667
$ws.Range("M132").Value = -2141.450000000001
$ws.Range("N132").Value = -7654.6667
$ws.Range("H135").Value = 370929.72
$ws.Range("I135").Value = 417033.4
$ws.Range("K135").Value = 3753300.6
$ws.Range("M135").Value = -3750765.6
$ws.Range("H138").Value = 4751.2856
$ws.Range("I138").Value = 883.75
$ws.Range("K138").Value = 2651.25
$ws.Range("M138").Value = 2488.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 731.5
$ws.Range("J16").Value = 500
$ws.Range("L16").Value = 500
$ws.Range("N16").Value = -1074
$ws.Range("H61").Value = 6094.9414
$ws.Range("I61").Value = 3301.5833
$ws.Range("K61").Value = 3301.5833
$ws.Range("M61").Value = -3089.5833
$ws.Range("H132").Value = 5851.08
$ws.Range("I132").Value = 4338.9443
$ws.Range("J132").Value = 9739.429
$ws.Range("K132").Value = 13016.8329
$ws.Range("L132").Value = 29218.287
$ws.Range("M132").Value = -10486.8329
$ws.Range("N132").Value = -34278.287
$ws.Range("H136").Value = 6094.9414
$ws.Range("I136").Value = 3301.5833
$ws.Range("K136").Value = 9904.749899999999
$ws.Range("M136").Value = -7354.749899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4783.523
$ws.Range("J134").Value = 10199.25
$ws.Range("L134").Value = 30597.75
$ws.Range("N134").Value = -35667.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 6974.222
$ws.Range("I132").Value = 1893.4445
$ws.Range("J132").Value = 12055
$ws.Range("K132").Value = 5680.333500000001
$ws.Range("L132").Value = 36165
$ws.Range("M132").Value = -3150.333500000001
$ws.Range("N132").Value = -41225

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 333339000
$ws.Range("J87").Value = 12000
$ws.Range("L87").Value = 36000
$ws.Range("N87").Value = -38496
$ws.Range("H90").Value = 333339000
$ws.Range("J90").Value = 12000
$ws.Range("L90").Value = 108000
$ws.Range("N90").Value = -120480
$ws.Range("H116").Value = 2399.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 63705.625
$ws.Range("I80").Value = 996.1
$ws.Range("J80").Value = 168221.5
$ws.Range("K80").Value = 996.1
$ws.Range("L80").Value = 168221.5
$ws.Range("M80").Value = 1.899999999999977
$ws.Range("N80").Value = -170217.5
$ws.Range("H83").Value = 63705.625
$ws.Range("I83").Value = 996.1
$ws.Range("J83").Value = 168221.5
$ws.Range("K83").Value = 4980.5
$ws.Range("L83").Value = 841107.5
$ws.Range("M83").Value = 11.5
$ws.Range("N83").Value = -851091.5
$ws.Range("H102").Value = 3765.1943
$ws.Range("I102").Value = 3438.2593
$ws.Range("K102").Value = 3438.2593
$ws.Range("M102").Value = -1816.2593
$ws.Range("H126").Value = 2943.2666
$ws.Range("I126").Value = 2909.5
$ws.Range("J126").Value = 2965.7778
$ws.Range("K126").Value = 8728.5
$ws.Range("L126").Value = 8897.3334
$ws.Range("M126").Value = -6258.5
$ws.Range("N126").Value = -13837.3334
$ws.Range("H132").Value = 2630.8484
$ws.Range("I132").Value = 2590.276
$ws.Range("K132").Value = 7770.828
$ws.Range("M132").Value = -5240.828

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2845.5557
$ws.Range("J22").Value = 3518.3333
$ws.Range("L22").Value = 3518.3333
$ws.Range("N22").Value = -4108.3333
$ws.Range("H27").Value = 2845.5557
$ws.Range("J27").Value = 3518.3333
$ws.Range("L27").Value = 3518.3333
$ws.Range("N27").Value = -3732.3333
$ws.Range("H32").Value = 7500
$ws.Range("I32").Value = 7500
$ws.Range("K32").Value = 7500
$ws.Range("M32").Value = -7183
$ws.Range("H82").Value = 1815.5555
$ws.Range("I82").Value = 1617.7
$ws.Range("J82").Value = 2062.875
$ws.Range("K82").Value = 1617.7
$ws.Range("L82").Value = 2062.875
$ws.Range("M82").Value = -1256.7
$ws.Range("N82").Value = -2784.875
$ws.Range("H85").Value = 1815.5555
$ws.Range("I85").Value = 1617.7
$ws.Range("J85").Value = 2062.875
$ws.Range("K85").Value = 1617.7
$ws.Range("L85").Value = 2062.875
$ws.Range("M85").Value = -369.7
$ws.Range("N85").Value = -4558.875
$ws.Range("H93").Value = 2893.077
$ws.Range("I93").Value = 2777.4
$ws.Range("K93").Value = 2777.4
$ws.Range("M93").Value = -1529.4
$ws.Range("H100").Value = 3881.6
$ws.Range("J100").Value = 6300.2
$ws.Range("L100").Value = 6300.2
$ws.Range("N100").Value = -7382.2
$ws.Range("H132").Value = 20008484
$ws.Range("I132").Value = 35718908
$ws.Range("J132").Value = 13394.363
$ws.Range("K132").Value = 107156724
$ws.Range("L132").Value = 40183.089
$ws.Range("M132").Value = -107154194
$ws.Range("N132").Value = -45243.089

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 27151.182
$ws.Range("I2").Value = 27151.182
$ws.Range("K2").Value = 27151.182
$ws.Range("M2").Value = -27039.182
$ws.Range("H43").Value = 15384.615
$ws.Range("J43").Value = 15384.615
$ws.Range("L43").Value = 15384.615
$ws.Range("N43").Value = -15682.615
$ws.Range("H122").Value = 4517.5
$ws.Range("I122").Value = 2404.4546
$ws.Range("J122").Value = 7100.1113
$ws.Range("K122").Value = 7213.3638
$ws.Range("L122").Value = 21300.3339
$ws.Range("M122").Value = -4763.3638
$ws.Range("N122").Value = -26200.3339
